$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.190.41"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.595.90"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.64"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.03"
$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.618.28"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  -0.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.101"
$ws.Range("E11").Value = "  -0.74%  "

$ws.Range("E12").Value = "  -2.49%  "

$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.055.16"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.150.82"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.48"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.609.23"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.03"
$ws.Range("E19").Value = "  +2.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.35"
$ws.Range("E20").Value = "  -1.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.33"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("E22").Value = "  +2.31%  "

$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.35"
$ws.Range("E24").Value = "  +3.99%  "

$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.403"
$ws.Range("E26").Value = "  -3.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.718.86"
$ws.Range("E27").Value = "  -0.67%  "

$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0750"
$ws.Range("E30").Value = "  -3.86%  "

$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("E32").Value = "  -5.03%  "

$ws.Range("E33").Value = "  +0.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.87"
$ws.Range("E34").Value = "  +1.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.76"
$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.04"
$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.868"
$ws.Range("E38").Value = "  -2.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.845"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("E40").Value = "  +2.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.10"
$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.56"
$ws.Range("E42").Value = "  -1.23%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.36%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.99"
$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0956"
$ws.Range("E46").Value = "  -0.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.63"
$ws.Range("E47").Value = "  +0.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.91"
$ws.Range("E48").Value = "  -0.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.10"
$ws.Range("E50").Value = "  +5.04%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.71"
$ws.Range("E51").Value = "  +2.05%  "
